# BIR libs configuration - no predictive model
#
# Adds a new "COUNTERPARTY_BIR" row to the "r CustomerUnit_AnalysisUnit"
# sheet (the relation table between Customer_Unit and Analysis_Unit),
# mirroring the existing BE-COUNTERPARTY_* rows, and leaves the workbook
# with the "Fields definition" sheet as the active tab / selection on
# F7 in the relation sheet (matching the last-saved UI state in the
# source workbook).

$wb = $excel.ActiveWorkbook

# "r CustomerUnit_AnalysisUnit" is the 2nd sheet.
$wsRelation = $wb.Worksheets.Item(2)

# New row 6: BE / COUNTERPARTY_BIR entry, following the same pattern as
# the existing BE-COUNTERPARTY_KOPER / BE-COUNTERPARTY_BIB rows.
$wsRelation.Range("A6").Value = "CREATE/MODIFY"
$wsRelation.Range("B6").Value = "BE-COUNTERPARTY_BIR"
$wsRelation.Range("C6").Value = "BE-COUNTERPARTY_BIR"
$wsRelation.Range("D6").Value = "BE-COUNTERPARTY_BIR"
$wsRelation.Range("E6").Value = "BE"
$wsRelation.Range("F6").Value = "COUNTERPARTY_BIR"

# Update the selection on this sheet to match the last-known cursor
# position (F7) before moving focus away from it.
$wsRelation.Range("F7").Select() | Out-Null

# "Fields definition" is the 3rd sheet; make it the active tab, as in
# the saved workbook (activeTab moves from index 1 to index 2).
$wsFields = $wb.Worksheets.Item(3)
$wsFields.Activate() | Out-Null
